$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two pairs of rows (20/21 and 58/59) had their data (columns B..AD)
# swapped between each other, while column A (the running index) stayed put.

$row20 = @{
    'B'  = 6228598
    'E'  = 'Geylang International'
    'F'  = 'Tampines Rovers FC'
    'G'  = 1
    'H'  = 1
    'I'  = 1
    'J'  = 1
    'K'  = 'D'
    'L'  = 3.6
    'M'  = 4.2
    'N'  = 1.666
    'O'  = 4
    'P'  = 4.5
    'Q'  = 1.55
    'R'  = 1
    'S'  = 1.85
    'T'  = 2
    'U'  = 3.25
    'V'  = 2.025
    'W'  = 1.825
    'X'  = -1
    'Y'  = 3.5
    'Z'  = -1
    'AA' = 0.8500000000000001
    'AB' = -1
    'AC' = -1
    'AD' = 0.825
}

$row21 = @{
    'B'  = 6228597
    'E'  = 'Hougang United FC'
    'F'  = 'Balestier Khalsa FC'
    'G'  = 1
    'H'  = 3
    'I'  = 0
    'J'  = 0
    'K'  = 'A'
    'L'  = 2.5
    'M'  = 3.6
    'N'  = 2.25
    'O'  = 2.6
    'P'  = 3.75
    'Q'  = 2.2
    'R'  = 0.25
    'S'  = 1.825
    'T'  = 2.025
    'U'  = 4
    'V'  = 1.95
    'W'  = 1.9
    'X'  = -1
    'Y'  = -1
    'Z'  = 1.2
    'AA' = -1
    'AB' = 1.025
    'AC' = 0
    'AD' = 0
}

$row58 = @{
    'B'  = 8089710
    'E'  = 'Balestier Khalsa FC'
    'F'  = 'Geylang International'
    'G'  = 2
    'H'  = 2
    'I'  = 0
    'J'  = 2
    'K'  = 'D'
    'L'  = 2.25
    'M'  = 4
    'N'  = 2.4
    'O'  = 2
    'P'  = 4
    'Q'  = 2.9
    'R'  = -0.5
    'S'  = 2.05
    'T'  = 1.8
    'U'  = 4
    'V'  = 1.825
    'W'  = 2.025
    'X'  = -1
    'Y'  = 3
    'Z'  = -1
    'AA' = -1
    'AB' = 0.8
    'AC' = 0
    'AD' = 0
}

$row59 = @{
    'B'  = 8088722
    'E'  = 'Hougang United FC'
    'F'  = 'Lion City Sailors FC'
    'G'  = 1
    'H'  = 4
    'I'  = 1
    'J'  = 2
    'K'  = 'A'
    'L'  = 7.5
    'M'  = 5.5
    'N'  = 1.25
    'O'  = 9.5
    'P'  = 7
    'Q'  = 1.181
    'R'  = 2.25
    'S'  = 1.875
    'T'  = 1.975
    'U'  = 4.25
    'V'  = 1.9
    'W'  = 1.95
    'X'  = -1
    'Y'  = -1
    'Z'  = 0.181
    'AA' = -1
    'AB' = 0.9750000000000001
    'AC' = 0.8999999999999999
    'AD' = -1
}

foreach ($col in $row20.Keys) {
    $ws.Range("$col" + "20").Value = $row20[$col]
}
foreach ($col in $row21.Keys) {
    $ws.Range("$col" + "21").Value = $row21[$col]
}
foreach ($col in $row58.Keys) {
    $ws.Range("$col" + "58").Value = $row58[$col]
}
foreach ($col in $row59.Keys) {
    $ws.Range("$col" + "59").Value = $row59[$col]
}
